# Updated symbol list on Wed Dec 21 07:42:49 UTC 2022 with GitHub Actions
#
# Refresh the crypto price/volume table. All "Price" (column D) cells are
# stored as literal text in the workbook (not numbers), so each one has its
# number format forced to Text ("@") before the new value is written -- this
# keeps things like trailing zeros ("3.700") and tiny decimals
# ("0.00005803") intact instead of having Excel re-parse them as floats
# (which would silently drop the trailing zero / switch to scientific
# notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "248.38"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.62"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.388"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05687"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.409"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.322"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8119"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9312"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1411"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07482"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03083"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03015"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09373"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.765"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001576"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04757"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0005852"
$ws.Range("E19").Value = "18OneONEWorstin24h"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006431"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004993"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.001024"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.700"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.3256"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1299"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03993"

# Row 41: was KickToken -> now BKEXToken
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1069"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42: was BKEXToken -> now CEJI
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002711"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43: was CEJI -> now KickToken
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.006872"
$ws.Range("E43").Value = "42KickTokenKICK"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007481"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005803"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4302"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2124"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
